# Update Case_4_108 vm_pu results: B (vm_pu at bus 0) changed from 1.05 to 1.02 (380 kV case),
# and all downstream bus voltages recomputed accordingly. Column A (row index), G, H unchanged.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$blockBF = New-Object "object[,]" 24,5
$blockBF[0,0] = 1.02
$blockBF[0,1] = 1.052067075342839
$blockBF[0,2] = 1.057668645511186
$blockBF[0,3] = 1.065070349420577
$blockBF[0,4] = 1.07108362620493
$blockBF[1,0] = 1.02
$blockBF[1,1] = 1.053157994356217
$blockBF[1,2] = 1.058520317407273
$blockBF[1,3] = 1.066101332288523
$blockBF[1,4] = 1.072088095464427
$blockBF[2,0] = 1.02
$blockBF[2,1] = 1.053864305389149
$blockBF[2,2] = 1.059071565936704
$blockBF[2,3] = 1.066769136319237
$blockBF[2,4] = 1.072738596733165
$blockBF[3,0] = 1.02
$blockBF[3,1] = 1.054161337533796
$blockBF[3,2] = 1.059303348378529
$blockBF[3,3] = 1.067050045592785
$blockBF[3,4] = 1.073012196719376
$blockBF[4,0] = 1.02
$blockBF[4,1] = 1.054211216342559
$blockBF[4,2] = 1.059342267852284
$blockBF[4,3] = 1.067097221124462
$blockBF[4,4] = 1.073058142922422
$blockBF[5,0] = 1.02
$blockBF[5,1] = 1.053868273958245
$blockBF[5,2] = 1.059074662878344
$blockBF[5,3] = 1.066772889197149
$blockBF[5,4] = 1.072742252082471
$blockBF[6,0] = 1.02
$blockBF[6,1] = 1.052435670913846
$blockBF[6,2] = 1.057956438487524
$blockBF[6,3] = 1.065418631889841
$blockBF[6,4] = 1.071422978594269
$blockBF[7,0] = 1.02
$blockBF[7,1] = 1.049914407888423
$blockBF[7,2] = 1.055987239765874
$blockBF[7,3] = 1.063037558923517
$blockBF[7,4] = 1.069102432123855
$blockBF[8,0] = 1.02
$blockBF[8,1] = 1.04823569171078
$blockBF[8,2] = 1.054675325027377
$blockBF[8,3] = 1.061453770395294
$blockBF[8,4] = 1.067558247371198
$blockBF[9,0] = 1.02
$blockBF[9,1] = 1.047509291665645
$blockBF[9,2] = 1.054107469150335
$blockBF[9,3] = 1.060768830212033
$blockBF[9,4] = 1.066890280007948
$blockBF[10,0] = 1.02
$blockBF[10,1] = 1.047239548486879
$blockBF[10,2] = 1.05389657459217
$blockBF[10,3] = 1.060514541093245
$blockBF[10,4] = 1.066642269024342
$blockBF[11,0] = 1.02
$blockBF[11,1] = 1.04729740596328
$blockBF[11,2] = 1.053941810715569
$blockBF[11,3] = 1.060569081167112
$blockBF[11,4] = 1.066695463606642
$blockBF[12,0] = 1.02
$blockBF[12,1] = 1.047486993099301
$blockBF[12,2] = 1.054090035876162
$blockBF[12,3] = 1.06074780797762
$blockBF[12,4] = 1.066869777262626
$blockBF[13,0] = 1.02
$blockBF[13,1] = 1.04760381380805
$blockBF[13,2] = 1.054181366599761
$blockBF[13,3] = 1.060857944457676
$blockBF[13,4] = 1.066977191167002
$blockBF[14,0] = 1.02
$blockBF[14,1] = 1.048283910900133
$blockBF[14,2] = 1.054713016263951
$blockBF[14,3] = 1.061499245572921
$blockBF[14,4] = 1.067602592421485
$blockBF[15,0] = 1.02
$blockBF[15,1] = 1.048710650224042
$blockBF[15,2] = 1.055046563175086
$blockBF[15,3] = 1.061901744825767
$blockBF[15,4] = 1.067995071093669
$blockBF[16,0] = 1.02
$blockBF[16,1] = 1.04895960784177
$blockBF[16,2] = 1.055241135848658
$blockBF[16,3] = 1.062136597886655
$blockBF[16,4] = 1.068224062343373
$blockBF[17,0] = 1.02
$blockBF[17,1] = 1.049044504069395
$blockBF[17,2] = 1.055307483545632
$blockBF[17,3] = 1.062216690618727
$blockBF[17,4] = 1.068302153466915
$blockBF[18,0] = 1.02
$blockBF[18,1] = 1.048664860191741
$blockBF[18,2] = 1.055010774633572
$blockBF[18,3] = 1.061858551975023
$blockBF[18,4] = 1.067952955111604
$blockBF[19,0] = 1.02
$blockBF[19,1] = 1.047431162368389
$blockBF[19,2] = 1.054046386361862
$blockBF[19,3] = 1.060695173845816
$blockBF[19,4] = 1.066818443420777
$blockBF[20,0] = 1.02
$blockBF[20,1] = 1.046655916419562
$blockBF[20,2] = 1.053440224907663
$blockBF[20,3] = 1.059964453525752
$blockBF[20,4] = 1.066105720310741
$blockBF[21,0] = 1.02
$blockBF[21,1] = 1.047066848397687
$blockBF[21,2] = 1.053761544579059
$blockBF[21,3] = 1.060351751771256
$blockBF[21,4] = 1.066483492295244
$blockBF[22,0] = 1.02
$blockBF[22,1] = 1.048685550588484
$blockBF[22,2] = 1.055026945871918
$blockBF[22,3] = 1.061878068711451
$blockBF[22,4] = 1.067971985309899
$blockBF[23,0] = 1.02
$blockBF[23,1] = 1.050565839904383
$blockBF[23,2] = 1.0564961719018
$blockBF[23,3] = 1.063652492123256
$blockBF[23,4] = 1.069701849873498
$ws.Range("B2:F25").Value = $blockBF

$blockIN = New-Object "object[,]" 24,6
$blockIN[0,0] = 1.042529645657527
$blockIN[0,1] = 1.057091632308688
$blockIN[0,2] = 1.060402984148173
$blockIN[0,3] = 1.067784580832454
$blockIN[0,4] = 1.073781743405864
$blockIN[0,5] = 1.058592823676671
$blockIN[1,0] = 1.042742864257942
$blockIN[1,1] = 1.057832198055621
$blockIN[1,2] = 1.061068441721175
$blockIN[1,3] = 1.068630351910447
$blockIN[1,4] = 1.074602233548445
$blockIN[1,5] = 1.059334441111909
$blockIN[2,0] = 1.042879581719978
$blockIN[2,1] = 1.058311221562279
$blockIN[2,2] = 1.061498551843274
$blockIN[2,3] = 1.069177700912126
$blockIN[2,4] = 1.075133070926037
$blockIN[2,5] = 1.059814144886891
$blockIN[3,0] = 1.042936758383796
$blockIN[3,1] = 1.058512561487079
$blockIN[3,2] = 1.061679253437947
$blockIN[3,3] = 1.069407824788559
$blockIN[3,4] = 1.075356216581179
$blockIN[3,5] = 1.060015770737478
$blockIN[4,0] = 1.04294634103442
$blockIN[4,1] = 1.058546364898805
$blockIN[4,2] = 1.061709587173263
$blockIN[4,3] = 1.069446464663361
$blockIN[4,4] = 1.075393682625873
$blockIN[4,5] = 1.060049622153927
$blockIN[5,0] = 1.042880346893986
$blockIN[5,1] = 1.05831391203996
$blockIN[5,2] = 1.061500966845992
$blockIN[5,3] = 1.069180775767062
$blockIN[5,4] = 1.075136052682064
$blockIN[5,5] = 1.059816839185359
$blockIN[6,0] = 1.042601962369024
$blockIN[6,1] = 1.057341945373075
$blockIN[6,2] = 1.060627978867062
$blockIN[6,3] = 1.068070396590329
$blockIN[6,4] = 1.074059047028346
$blockIN[6,5] = 1.058843492214321
$blockIN[7,0] = 1.042101858348756
$blockIN[7,1] = 1.055627914218465
$blockIN[7,2] = 1.059085966097008
$blockIN[7,3] = 1.06611438795261
$blockIN[7,4] = 1.072160676622867
$blockIN[7,5] = 1.057127026938869
$blockIN[8,0] = 1.041762047848088
$blockIN[8,1] = 1.054484371335485
$blockIN[8,2] = 1.058055499498231
$blockIN[8,3] = 1.064810824505765
$blockIN[8,4] = 1.070894755897363
$blockIN[8,5] = 1.055981860093829
$blockIN[9,0] = 1.041613389492006
$blockIN[9,1] = 1.053989005299041
$blockIN[9,2] = 1.057608718248073
$blockIN[9,3] = 1.064246476585429
$blockIN[9,4] = 1.070346522895357
$blockIN[9,5] = 1.055485790580795
$blockIN[10,0] = 1.041557943282508
$blockIN[10,1] = 1.053804973793456
$blockIN[10,2] = 1.057442676642891
$blockIN[10,3] = 1.064036868498064
$blockIN[10,4] = 1.070142872939022
$blockIN[10,5] = 1.055301497729363
$blockIN[11,0] = 1.041569846991182
$blockIN[11,1] = 1.053844450553089
$blockIN[11,2] = 1.057478297069279
$blockIN[11,3] = 1.064081829422137
$blockIN[11,4] = 1.070186557074621
$blockIN[11,5] = 1.055341030550523
$blockIN[12,0] = 1.041608810938431
$blockIN[12,1] = 1.053973793818427
$blockIN[12,2] = 1.057594994973145
$blockIN[12,3] = 1.064229149990388
$blockIN[12,4] = 1.070329689365277
$blockIN[12,5] = 1.055470557498133
$blockIN[13,0] = 1.04163278772271
$blockIN[13,1] = 1.054053482437342
$blockIN[13,2] = 1.057666884864576
$blockIN[13,3] = 1.06431992116456
$blockIN[13,4] = 1.070417876339351
$blockIN[13,5] = 1.055550359284027
$blockIN[14,0] = 1.041771881845471
$blockIN[14,1] = 1.054517242858697
$blockIN[14,2] = 1.05808513865701
$blockIN[14,3] = 1.064848280626753
$blockIN[14,4] = 1.070931138640468
$blockIN[14,5] = 1.056014778298375
$blockIN[15,0] = 1.041858725605502
$blockIN[15,1] = 1.054808092935822
$blockIN[15,2] = 1.058347342525642
$blockIN[15,3] = 1.065179734510175
$blockIN[15,4] = 1.071253073093656
$blockIN[15,5] = 1.05630604141597
$blockIN[16,0] = 1.041909233606154
$blockIN[16,1] = 1.054977721127166
$blockIN[16,2] = 1.058500225388462
$blockIN[16,3] = 1.065373075934544
$blockIN[16,4] = 1.071440844176698
$blockIN[16,5] = 1.056475910498802
$blockIN[17,0] = 1.041926430674616
$blockIN[17,1] = 1.055035556602633
$blockIN[17,2] = 1.058552344947882
$blockIN[17,3] = 1.065439002048333
$blockIN[17,4] = 1.071504867923511
$blockIN[17,5] = 1.056533828107279
$blockIN[18,0] = 1.041849423243079
$blockIN[18,1] = 1.054776889515119
$blockIN[18,2] = 1.058319216344301
$blockIN[18,3] = 1.065144171610716
$blockIN[18,4] = 1.071218533381533
$blockIN[18,5] = 1.056274793682831
$blockIN[19,0] = 1.041597343315941
$blockIN[19,1] = 1.053935706290942
$blockIN[19,2] = 1.057560632752228
$blockIN[19,3] = 1.064185767314855
$blockIN[19,4] = 1.070287540796284
$blockIN[19,5] = 1.05543241588199
$blockIN[20,0] = 1.041437532339907
$blockIN[20,1] = 1.053406644632623
$blockIN[20,2] = 1.057083177054381
$blockIN[20,3] = 1.063583272386337
$blockIN[20,4] = 1.069702120659625
$blockIN[20,5] = 1.054902602895432
$blockIN[21,0] = 1.041522376052914
$blockIN[21,1] = 1.053687126876361
$blockIN[21,2] = 1.057336333019467
$blockIN[21,3] = 1.063902657541622
$blockIN[21,4] = 1.070012469290782
$blockIN[21,5] = 1.055183483456128
$blockIN[22,0] = 1.04185362703287
$blockIN[22,1] = 1.054790989058712
$blockIN[22,2] = 1.058331925528627
$blockIN[22,3] = 1.065160240923283
$blockIN[22,4] = 1.071234140414854
$blockIN[22,5] = 1.056288913249393
$blockIN[23,0] = 1.042232277851301
$blockIN[23,1] = 1.05607118458682
$blockIN[23,2] = 1.059485048607568
$blockIN[23,3] = 1.066619986914063
$blockIN[23,4] = 1.072651512979423
$blockIN[23,5] = 1.057570926801991
$ws.Range("I2:N25").Value = $blockIN

Write-Output "Updated vm_pu results for case with 380 kV (rows 2-25)"
